$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.275.64'
$ws.Cells.Item(2, 5).Value = '  +0.16%  '
$ws.Cells.Item(3, 4).Value = '1.904.94'
$ws.Cells.Item(3, 5).Value = '  +1.76%  '
$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).Value = '314.46'
$ws.Cells.Item(5, 5).Value = '  +0.80%  '
$ws.Cells.Item(6, 4).Value = '1.002'
$ws.Cells.Item(6, 5).Value = '  +0.04%  '
$ws.Cells.Item(7, 4).Value = '0.5078'
$ws.Cells.Item(7, 5).Value = '  +1.26%  '
$ws.Cells.Item(8, 4).Value = '0.3936'
$ws.Cells.Item(8, 5).Value = '  -0.48%  '
$ws.Cells.Item(9, 4).Value = '0.09655'
$ws.Cells.Item(9, 5).Value = '  -2.18%  '
$ws.Cells.Item(10, 5).Value = '  +0.21%  '
$ws.Cells.Item(11, 4).Value = '42.11'
$ws.Cells.Item(11, 5).Value = '  +2.03%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '6.430'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.67%  '
$ws.Cells.Item(13, 2).Value = 'Solana'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(13, 4).Value = '20.97'
$ws.Cells.Item(13, 5).Value = '  -0.16%  '
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value = '1.917.70'
$ws.Cells.Item(14, 5).Value = '  +2.63%  '
$ws.Cells.Item(15, 4).Value = '7.331'
$ws.Cells.Item(15, 5).Value = '  -0.95%  '
$ws.Cells.Item(16, 4).Value = '1.005'
$ws.Cells.Item(16, 5).Value = '  +0.26%  '
$ws.Cells.Item(17, 5).Value = '  -1.15%  '
$ws.Cells.Item(18, 4).Value = '92.82'
$ws.Cells.Item(18, 5).Value = '  -0.84%  '
$ws.Cells.Item(19, 4).Value = '0.06636'
$ws.Cells.Item(19, 5).Value = '  -0.11%  '
$ws.Cells.Item(20, 4).Value = '18.04'
$ws.Cells.Item(20, 5).Value = '  +3.34%  '
$ws.Cells.Item(21, 4).Value = '1.002'
$ws.Cells.Item(21, 5).Value = '  +0.07%  '
$ws.Cells.Item(22, 4).Value = '6.232'
$ws.Cells.Item(22, 5).Value = '  +1.88%  '
$ws.Cells.Item(23, 4).Value = '28.350.31'
$ws.Cells.Item(23, 5).Value = '  +0.24%  '
$ws.Cells.Item(24, 4).Value = '11.33'
$ws.Cells.Item(24, 5).Value = '  -0.16%  '
$ws.Cells.Item(25, 5).Value = '  +1.86%  '
$ws.Cells.Item(26, 2).Value = 'LEO'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(26, 4).Value = '3.382'
$ws.Cells.Item(26, 5).Value = '  -1.41%  '
$ws.Cells.Item(27, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(27, 4).Value = '2.664'
$ws.Cells.Item(27, 5).Value = '  +4.07%  '
$ws.Cells.Item(28, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(28, 4).Value = '2.141.20'
$ws.Cells.Item(28, 5).Value = '  +2.80%  '
$ws.Cells.Item(29, 2).Value = 'EthereumClassic'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(29, 4).Value = '21.02'
$ws.Cells.Item(29, 5).Value = '  -1.15%  '
$ws.Cells.Item(30, 2).Value = 'Monero'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(30, 4).Value = '158.12'
$ws.Cells.Item(30, 5).Value = '  -0.08%  '
$ws.Cells.Item(31, 2).Value = 'BitcoinCash'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(31, 4).Value = '127.01'
$ws.Cells.Item(31, 5).Value = '  -0.79%  '
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).Value = '1.099'
$ws.Cells.Item(32, 5).Value = '  +3.69%  '
$ws.Cells.Item(33, 2).Value = 'Stellar'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(33, 4).Value = '0.1066'
$ws.Cells.Item(33, 5).Value = '  +0.56%  '
$ws.Cells.Item(34, 2).Value = 'Filecoin'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(34, 4).Value = '5.655'
$ws.Cells.Item(34, 5).Value = '  +0.34%  '
$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(35, 4).Value = '3.631'
$ws.Cells.Item(35, 5).Value = '  +0.64%  '
$ws.Cells.Item(36, 2).Value = 'FraxShare'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(36, 4).Value = '9.683'
$ws.Cells.Item(36, 5).Value = '  +2.09%  '
$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(37, 4).Value = '0.06676'
$ws.Cells.Item(37, 5).Value = '  -1.97%  '
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(38, 4).Value = '0.02427'
$ws.Cells.Item(38, 5).Value = '  +0.99%  '
$ws.Cells.Item(39, 2).Value = 'ARBITRUM'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(39, 4).Value = '1.245'
$ws.Cells.Item(39, 5).Value = '  +2.35%  '
$ws.Cells.Item(40, 2).Value = 'Algorand'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(40, 4).Value = '0.2204'
$ws.Cells.Item(40, 5).Value = '  +0.75%  '
$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(41, 4).Value = '1.283'
$ws.Cells.Item(41, 5).Value = '  +9.17%  '
$ws.Cells.Item(42, 2).Value = 'TheSandbox'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(42, 4).Value = '0.6395'
$ws.Cells.Item(42, 5).Value = '  +1.40%  '
$ws.Cells.Item(43, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(43, 4).Value = '5.022'
$ws.Cells.Item(43, 5).Value = '  +0.11%  '
$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(44, 4).Value = '11.51'
$ws.Cells.Item(44, 5).Value = '  -0.04%  '
$ws.Cells.Item(45, 2).Value = 'Frax'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(45, 4).Value = '1.002'
$ws.Cells.Item(45, 5).Value = '  +0.13%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).Value = '13.55'
$ws.Cells.Item(46, 5).Value = '  +0.41%  '
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = '0.6046'
$ws.Cells.Item(47, 5).Value = '  +0.45%  '
$ws.Cells.Item(48, 2).Value = 'PancakeSwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(48, 4).Value = '3.754'
$ws.Cells.Item(48, 5).Value = '  +2.43%  '
$ws.Cells.Item(49, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(49, 4).Value = '1.284'
$ws.Cells.Item(49, 5).Value = '  +0.98%  '
$ws.Cells.Item(50, 2).Value = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(50, 4).Value = '2.048'
$ws.Cells.Item(50, 5).Value = '  +2.80%  '
$ws.Cells.Item(51, 2).Value = 'Quant'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(51, 4).Value = '123.65'
$ws.Cells.Item(51, 5).Value = '  -1.00%  '
